$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 9
$ws.Range("F4").Value = 3460
$ws.Range("F5").Value = 3460
$ws.Range("F6").Value = 240
$ws.Range("F7").Value = 4983
$ws.Range("F8").Value = 500
$ws.Range("F10").Value = 189
$ws.Range("F11").Value = 663
$ws.Range("F16").Value = 297
$ws.Range("F17").Value = 30
$ws.Range("F19").Value = 154
$ws.Range("F21").Value = 355
$ws.Range("F22").Value = 4835
$ws.Range("F23").Value = 39
$ws.Range("F26").Value = 5966
$ws.Range("F28").Value = 14
$ws.Range("F29").Value = 3209
$ws.Range("F30").Value = 304
$ws.Range("F31").Value = 693
$ws.Range("F32").Value = 4436
$ws.Range("F33").Value = 313
$ws.Range("F34").Value = 109
$ws.Range("F36").Value = 938
$ws.Range("F37").Value = 79
$ws.Range("F38").Value = 20
$ws.Range("F40").Value = 837
$ws.Range("F41").Value = 921
$ws.Range("F42").Value = 13

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 42
$ws.Range("F4").Value = 20
$ws.Range("F5").Value = 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1104
$ws.Range("F4").Value = 46

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1104
$ws.Range("F5").Value = 46
$ws.Range("F6").Value = 9
$ws.Range("F8").Value = 3460
$ws.Range("F9").Value = 3460
$ws.Range("F10").Value = 240
$ws.Range("F11").Value = 4983
$ws.Range("F12").Value = 500
$ws.Range("F14").Value = 189
$ws.Range("F15").Value = 663
$ws.Range("F19").Value = 297
$ws.Range("F20").Value = 30
$ws.Range("F21").Value = 42
$ws.Range("F23").Value = 154
$ws.Range("F25").Value = 355
$ws.Range("F26").Value = 4836
$ws.Range("F27").Value = 39
$ws.Range("F30").Value = 5966
$ws.Range("F32").Value = 14
$ws.Range("F33").Value = 3209
$ws.Range("F34").Value = 304
$ws.Range("F35").Value = 693
$ws.Range("F36").Value = 4436
$ws.Range("F37").Value = 313
$ws.Range("F38").Value = 20
$ws.Range("F39").Value = 109
$ws.Range("F41").Value = 938
$ws.Range("F42").Value = 79
$ws.Range("F43").Value = 20
$ws.Range("F45").Value = 837
$ws.Range("F46").Value = 921
$ws.Range("F47").Value = 11
$ws.Range("F48").Value = 13
